$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.170.52"
$ws.Range("E2").Value = "  +2.60%  "

# Row 3
$ws.Range("D3").Value = "2.301.71"
$ws.Range("E3").Value = "  +2.27%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$c = $ws.Range("D5")
$c.Value = "'308.13"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "

# Row 6
$c = $ws.Range("D6")
$c.Value = "'104.22"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +7.97%  "

# Row 7
$c = $ws.Range("D7")
$c.Value = "'0.527"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.55%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$c = $ws.Range("D9")
$c.Value = "'0.520"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +7.00%  "

# Row 10
$c = $ws.Range("D10")
$c.Value = "'35.88"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.69%  "

# Row 11
$c = $ws.Range("D11")
$c.Value = "'52.17"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.26%  "

# Row 12
$c = $ws.Range("D12")
$c.Value = "'0.0807"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "

# Row 13
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$c = $ws.Range("D14")
$c.Value = "'6.93"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.01%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.607.14"
$ws.Range("E15").Value = "  +15.70%  "

# Row 16
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.659.06"
$ws.Range("E16").Value = "  +2.27%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D17")
$c.Value = "'15.02"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.75%  "

# Row 18
$c = $ws.Range("D18")
$c.Value = "'0.799"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.59%  "

# Row 19
$ws.Range("D19").Value = "43.135.18"
$ws.Range("E19").Value = "  +2.84%  "

# Row 20
$c = $ws.Range("D20")
$c.Value = "'11.85"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.71%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0919"
$ws.Range("E21").Value = "  +1.71%  "

# Row 22
$c = $ws.Range("D22")
$c.Value = "'6.15"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.75%  "

# Row 23
$c = $ws.Range("D23")
$c.Value = "'67.67"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.76%  "

# Row 24
$c = $ws.Range("D24")
$c.Value = "'240.12"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.26%  "

# Row 25
$c = $ws.Range("D25")
$c.Value = "'2.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.43%  "

# Row 26
$c = $ws.Range("D26")
$c.Value = "'2.59"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.77%  "

# Row 27
$ws.Range("E27").Value = "  +0.46%  "

# Row 28
$c = $ws.Range("D28")
$c.Value = "'24.68"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +5.64%  "

# Row 29
$ws.Range("E29").Value = "  +5.47%  "

# Row 30
$c = $ws.Range("D30")
$c.Value = "'35.99"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.81%  "

# Row 31
$c = $ws.Range("D31")
$c.Value = "'9.53"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "

# Row 32
$c = $ws.Range("D32")
$c.Value = "'161.40"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.96%  "

# Row 33
$c = $ws.Range("D33")
$c.Value = "'5.20"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.47%  "

# Row 34
$c = $ws.Range("D34")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "

# Row 35
$c = $ws.Range("D35")
$c.Value = "'18.14"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.87%  "

# Row 36
$ws.Range("E36").Value = "  +7.06%  "

# Row 37
$c = $ws.Range("D37")
$c.Value = "'0.0733"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.54%  "

# Row 38
$c = $ws.Range("D38")
$c.Value = "'2.99"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.78%  "

# Row 39
$c = $ws.Range("D39")
$c.Value = "'4.51"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +10.35%  "

# Row 40
$ws.Range("E40").Value = "  +2.73%  "

# Row 41
$c = $ws.Range("D41")
$c.Value = "'1.85"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.50%  "

# Row 42
$ws.Range("E42").Value = "  +0.12%  "

# Row 43
$c = $ws.Range("D43")
$c.Value = "'2.46"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +11.15%  "

# Row 44
$ws.Range("D44").Value = "1.964.58"
$ws.Range("E44").Value = "  +1.36%  "

# Row 45
$c = $ws.Range("D45")
$c.Value = "'0.0288"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.14%  "

# Row 46
$c = $ws.Range("D46")
$c.Value = "'18.74"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

# Row 47
$c = $ws.Range("D47")
$c.Value = "'3.06"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +5.45%  "

# Row 48
$c = $ws.Range("D48")
$c.Value = "'10.18"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +5.54%  "

# Row 49
$c = $ws.Range("D49")
$c.Value = "'56.86"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.88%  "

# Row 50
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D50")
$c.Value = "'2.89"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.69%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D51")
$c.Value = "'1.57"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +7.59%  "
